# Merge in the "HashTag Time Series" sheet from the ICWSM paper results.

$wb = $excel.ActiveWorkbook

# Add a new worksheet after the last existing sheet ("Pred. vs Time").
$lastIndex = $wb.Worksheets.Count
$ws = $wb.Worksheets.Add($null, $wb.Worksheets.Item($lastIndex))
$ws.Name = "HashTag Time Series"

# Row labels in column A, written in the original authoring order so the
# shared-string table gets populated hashTag, beatles, facebook,
# elmundoconchavez, chavistas, ucvistasconchavez, vivachavez (note rows 7
# and 8 reuse the "vivachavez" label).
$ws.Cells.Item(1, 1).Value = "hashTag"
$ws.Cells.Item(3, 1).Value = "beatles"
$ws.Cells.Item(5, 1).Value = "facebook"
$ws.Cells.Item(2, 1).Value = "elmundoconchavez"
$ws.Cells.Item(4, 1).Value = "chavistas"
$ws.Cells.Item(6, 1).Value = "ucvistasconchavez"
$ws.Cells.Item(7, 1).Value = "vivachavez"
$ws.Cells.Item(8, 1).Value = "vivachavez"

# Days-to-election bucket header, followed by the per-hashtag time series
# (column B onward; column A already filled in above).
$data = @(
    @(24, 18, 15, 12, 9, 6, 3, 0),
    @(0.86956522802900005, 0.86956522234, 0.86956521974199996, 0.86956521850699997, 0.76512312000000005, 0.81234119999999999, 0.75123412000000001, 0.79123421299999996),
    @(0, 0.74474498239700004, 0.36763937841900002, 0, 0, 0, 0, 0),
    @(0.80032520618000003, 0.69593496157800006, 0.60516083007699994, 0.52622680608500005, 0.45758852640499997, 0.56098323999999999, 0.61324234, 0.55231229999999998),
    @(0.773297341995, 0, 0.78592625348599998, 0, 0, 0.37912323999999997, 0, 0),
    @(0, 0, 0.60129473282300006, 0.52286498240299994, 0.45466520146799999, 0.36763937800000002, 0.40233999999999998, 0.38912311999999999),
    @(0.62675629742500005, 0.43595703799300001, 0.34906271389499999, 0.52027991832300002, 0.45241874975800001, 0.55231229999999998, 0.69593496157800006, 0.36763937841900002),
    @(0, 0, 0.36323450000000002, 0.45241874975800001, 0.56098323999999999, 0.31323322999999997, 0.61234124320000005, 0.49912312399999997)
)

for ($r = 0; $r -lt $data.Length; $r++) {
    $row = $data[$r]
    for ($c = 0; $c -lt $row.Length; $c++) {
        $ws.Cells.Item($r + 1, $c + 2).Value = $row[$c]
    }
}

# Single standout cell (G6, row 6 / col 7) gets a plain-black font color override.
$ws.Cells.Item(6, 7).Font.Color = 0

# Column A is sized to fit the longest hashtag label ("ucvistasconchavez").
$ws.Columns("A").ColumnWidth = 19.43

# Match the author's last selection/cursor position on the new sheet.
$ws.Range("B9").Select() | Out-Null

# Restore the window geometry recorded by the author's Excel session.
$win = $excel.ActiveWindow
$win.Left = -80
$win.Top = -440
$win.Width = 25600
$win.Height = 16000
